# WIMR_Config_testing.xlsx - event code adjustments, rank of ICA, epoch
# binning and visualization.
#
# The underlying change is: the DownSample value (column D, "DownSample")
# for the SUB 810 config row is cleared out, and the active selection is
# moved from F9 to D4 (the now-empty DownSample cell) to reflect where the
# user was working/reviewing the config.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the DownSample value in D2 (was 256)
$ws.Range("D2").ClearContents()

# Reflect the user's new selection/viewport on the DownSample column
$ws.Range("D4").Select() | Out-Null
